$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 0.7333333333333333
$ws.Range("D3").Value = 0.8666666666666667
$ws.Range("F3").Value = 0.9333333333333333
$ws.Range("H3").Value = 0.5283018867924528
$ws.Range("I3").Value = 0.1214298291952809
$ws.Range("J3").Value = 0.6666666666666666
$ws.Range("K3").Value = 237.8666666666667

$ws.Range("Q3").Value = 62
$ws.Range("R3").Value = 62
$ws.Range("S3").Value = 87
$ws.Range("T3").Value = 113
$ws.Range("U3").Value = 210
$ws.Range("V3").Value = 1831
$ws.Range("W3").Value = 1831
$ws.Range("X3").Value = 1806
$ws.Range("Y3").Value = 1780
$ws.Range("Z3").Value = 1683

$ws.Range("AF3").Value = 0.967248
$ws.Range("AG3").Value = 0.967248
$ws.Range("AH3").Value = 0.954041
$ws.Range("AI3").Value = 0.940306
$ws.Range("AJ3").Value = 0.889065
